$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.092.06'
$ws.Range('E2').Value = '  -3.72%  '

$ws.Range('D3').Value = '3.508.64'
$ws.Range('E3').Value = '  -4.94%  '

$ws.Range('E4').Value = '  +0.02%  '

$cell = $ws.Range('D5')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '580.57'
$cell.Style = $origStyle
$ws.Range('E5').Value = '  -1.50%  '

$cell = $ws.Range('D6')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '174.74'
$cell.Style = $origStyle
$ws.Range('E6').Value = '  -2.97%  '

$cell = $ws.Range('D7')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.627'
$cell.Style = $origStyle
$ws.Range('E7').Value = '  +0.79%  '

$ws.Range('D8').Value = '3.501.80'
$ws.Range('E8').Value = '  -4.93%  '

$ws.Range('E9').Value = '  +0.07%  '

$ws.Range('E10').Value = '  -5.98%  '

$cell = $ws.Range('D11')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '6.72'
$cell.Style = $origStyle
$ws.Range('E11').Value = '  +6.26%  '

$cell = $ws.Range('D12')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.598'
$cell.Style = $origStyle
$ws.Range('E12').Value = '  -2.69%  '

$cell = $ws.Range('D13')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '47.26'
$cell.Style = $origStyle
$ws.Range('E13').Value = '  -5.58%  '

$ws.Range('E14').Value = '  -3.56%  '

$cell = $ws.Range('D15')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '677.74'
$cell.Style = $origStyle
$ws.Range('E15').Value = '  -1.21%  '

$ws.Range('D16').Value = '4.075.61'
$ws.Range('E16').Value = '  -4.90%  '

$cell = $ws.Range('D17')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '8.75'
$cell.Style = $origStyle
$ws.Range('E17').Value = '  -2.84%  '

$ws.Range('D18').Value = '3.524.53'
$ws.Range('E18').Value = '  -4.67%  '

$ws.Range('D19').Value = '69.104.92'
$ws.Range('E19').Value = '  -3.77%  '

$ws.Range('E20').Value = '  -1.41%  '

$cell = $ws.Range('D21')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '17.57'
$cell.Style = $origStyle
$ws.Range('E21').Value = '  -3.10%  '

$cell = $ws.Range('D22')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '11.22'
$cell.Style = $origStyle
$ws.Range('E22').Value = '  -4.00%  '

$cell = $ws.Range('D23')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.906'
$cell.Style = $origStyle
$ws.Range('E23').Value = '  -3.93%  '

$cell = $ws.Range('D24')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '16.06'
$cell.Style = $origStyle
$ws.Range('E24').Value = '  -10.20%  '

$cell = $ws.Range('D25')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '98.29'
$cell.Style = $origStyle
$ws.Range('E25').Value = '  -5.63%  '

$ws.Range('E26').Value = '  -4.19%  '

$cell = $ws.Range('D27')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '5.87'
$cell.Style = $origStyle
$ws.Range('E27').Value = '  +0.07%  '

$ws.Range('B28').Value = 'Dai'
$ws.Range('C28').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$cell = $ws.Range('D28')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = $origStyle
$ws.Range('E28').Value = '  -0.01%  '

$ws.Range('B29').Value = 'ImmutableX'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$cell = $ws.Range('D29')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.67'
$cell.Style = $origStyle
$ws.Range('E29').Value = '  -6.40%  '

$cell = $ws.Range('D30')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '9.49'
$cell.Style = $origStyle

$cell = $ws.Range('D31')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '32.93'
$cell.Style = $origStyle
$ws.Range('E31').Value = '  -7.43%  '

$cell = $ws.Range('D32')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '8.77'
$cell.Style = $origStyle
$ws.Range('E32').Value = '  -5.24%  '

$ws.Range('E33').Value = '  -7.67%  '

$cell = $ws.Range('D34')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '7.33'
$cell.Style = $origStyle
$ws.Range('E34').Value = '  -0.22%  '

$ws.Range('E35').Value = '  -5.51%  '

$cell = $ws.Range('D36')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '586.96'
$cell.Style = $origStyle
$ws.Range('E36').Value = '  +2.95%  '

$ws.Range('E37').Value = '  -15.28%  '

$cell = $ws.Range('D38')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '10.93'
$cell.Style = $origStyle
$ws.Range('E38').Value = '  -3.55%  '

$ws.Range('E39').Value = '  -4.27%  '

$cell = $ws.Range('D40')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '57.31'
$cell.Style = $origStyle
$ws.Range('E40').Value = '  -3.84%  '

$cell = $ws.Range('D41')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '1.00'
$cell.Style = $origStyle
$ws.Range('E41').Value = '  +0.15%  '

$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell = $ws.Range('D42')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.0442'
$cell.Style = $origStyle
$ws.Range('E42').Value = '  -4.92%  '

$ws.Range('B43').Value = 'TheGraph'
$ws.Range('C43').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$cell = $ws.Range('D43')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '0.339'
$cell.Style = $origStyle
$ws.Range('E43').Value = '  -4.08%  '

$ws.Range('E44').Value = '  -6.97%  '

$ws.Range('D45').Value = '3.435.17'
$ws.Range('E45').Value = '  -9.83%  '

$cell = $ws.Range('D46')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '33.48'
$cell.Style = $origStyle
$ws.Range('E46').Value = '  -5.66%  '

$ws.Range('D47').Value = '0.0₃0708'
$ws.Range('E47').Value = '  -9.57%  '

$cell = $ws.Range('D48')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.93'
$cell.Style = $origStyle
$ws.Range('E48').Value = '  +1.30%  '

$cell = $ws.Range('D49')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '2.60'
$cell.Style = $origStyle
$ws.Range('E49').Value = '  -7.16%  '

$ws.Range('E50').Value = '  -0.20%  '

$cell = $ws.Range('D51')
$origStyle = $cell.Style
$cell.NumberFormat = '@'
$cell.Value = '132.35'
$cell.Style = $origStyle
$ws.Range('E51').Value = '  -1.97%  '
